# The deck ships two themes:
#   ppt/theme/theme1.xml -> bound to the (only) Slide Master, name "Integral"
#   ppt/theme/theme2.xml -> bound to the Notes Master,        name "Office Theme"
#
# The authored change swaps the two themes' contents: the Slide Master's
# theme becomes the stock "Office Theme" color scheme, and the Notes
# Master's theme becomes the old "Integral" color scheme.
#
# The theme color palette is the only thing that differs between the two
# themes (font scheme and format scheme are byte-identical), so the edit
# is expressed purely as a 12-slot ThemeColorScheme.Colors(i).RGB swap on
# the Slide Master's theme - the documented way to edit a theme's colors
# in this host (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink, in that
# order). RGB values use the standard VBA COLORREF encoding
# (R + G*256 + B*65536).

$p = $ppt.ActivePresentation

# Target palette: stock PowerPoint "Office Theme" colors, which is what
# theme1.xml (the Slide Master's theme) must become.
$officeThemeColors = @(
    0,         # dk1      000000
    16777215,  # lt1      FFFFFF
    6968388,   # dk2      44546A
    15132391,  # lt2      E7E6E6
    13998939,  # accent1  5B9BD5
    3243501,   # accent2  ED7D31
    10855845,  # accent3  A5A5A5
    49407,     # accent4  FFC000
    12874308,  # accent5  4472C4
    4697456,   # accent6  70AD47
    12673797,  # hlink    0563C1
    7491477    # folHlink 954F72
)

$master = $p.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

for ($i = 1; $i -le 12; $i++) {
    $colorScheme.Colors($i).RGB = $officeThemeColors[$i - 1]
}
